# Estadisticos Segundo Parcial 26 Mayo
# Update the "2o Parcial" and "Final" sheets with the new statistics
# for rows 3 (2ALCV), 8 (Totales Docente - Ameca Garcia Ivan) and
# 29 (Totales Generales).

$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial" ---
$ws = $wb.Worksheets.Item("2o Parcial")

# Row 3
$ws.Range("E3").Value = 22
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 78.59999999999999
$ws.Range("H3").Value = 21.4
$ws.Range("I3").Value = 7.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

# Row 8
$ws.Range("E8").Value = 117
$ws.Range("F8").Value = 30
$ws.Range("G8").Value = 79.59999999999999
$ws.Range("H8").Value = 20.4
$ws.Range("I8").Value = 7.4
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# Row 29
$ws.Range("E29").Value = 554
$ws.Range("F29").Value = 80
$ws.Range("G29").Value = 87.40000000000001
$ws.Range("H29").Value = 12.6
$ws.Range("I29").Value = 8
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0

# --- Sheet "Final" ---
$ws = $wb.Worksheets.Item("Final")

# Row 3
$ws.Range("E3").Value = 22
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 78.59999999999999
$ws.Range("H3").Value = 21.4
$ws.Range("I3").Value = 6.8

# Row 8
$ws.Range("E8").Value = 117
$ws.Range("F8").Value = 30
$ws.Range("G8").Value = 79.59999999999999
$ws.Range("H8").Value = 20.4
$ws.Range("I8").Value = 6.9

# Row 29
$ws.Range("E29").Value = 554
$ws.Range("F29").Value = 80
$ws.Range("G29").Value = 87.40000000000001
$ws.Range("H29").Value = 12.6
$ws.Range("I29").Value = 7.8
